# Update column F ("dSF") values to reflect the re-pulled / recalculated data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    7  = 4
    9  = 1
    13 = -4
    15 = 3
    16 = 5
    21 = 12
    22 = 0
    24 = -3
    26 = -2
    30 = 5
    31 = 3
    35 = -2
    37 = 1
    39 = -2
    40 = -4
    45 = -3
    49 = 0
    55 = -2
    56 = 0
    59 = -3
    61 = -5
    63 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
